$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.816.45'
$ws.Range("E2").Value = '  +0.27%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.813.56'
$ws.Range("E3").Value = '  +0.78%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.25'
$ws.Range("E5").Value = '  +1.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.17'
$ws.Range("E6").Value = '  -0.42%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  +0.47%  '

$ws.Range("E10").Value = '  +1.19%  '

$ws.Range("E11").Value = '  +0.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000251'
$ws.Range("E12").Value = '  -0.57%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.05'
$ws.Range("E13").Value = '  +0.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.456.28'
$ws.Range("E14").Value = '  +0.82%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.818.68'
$ws.Range("E15").Value = '  +1.36%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.853.32'
$ws.Range("E16").Value = '  +0.31%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.40'
$ws.Range("E17").Value = '  -0.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.10'
$ws.Range("E18").Value = '  +1.24%  '

$ws.Range("E19").Value = '  +1.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '464.13'
$ws.Range("E20").Value = '  +0.97%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.86'
$ws.Range("E21").Value = '  -1.33%  '

$ws.Range("E22").Value = '  +1.13%  '

$ws.Range("E23").Value = '  -3.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.43'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.17'
$ws.Range("E25").Value = '  +1.63%  '

$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.04'
$ws.Range("E27").Value = '  +0.35%  '

$ws.Range("E28").Value = '  -0.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.964.87'
$ws.Range("E29").Value = '  +0.86%  '

$ws.Range("E30").Value = '  +0.76%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.41'
$ws.Range("E31").Value = '  +2.95%  '

$ws.Range("E32").Value = '  -0.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.53'
$ws.Range("E33").Value = '  -0.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.09'
$ws.Range("E35").Value = '  +0.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.1000'
$ws.Range("E36").Value = '  -0.15%  '

$ws.Range("E37").Value = '  +0.53%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.82'
$ws.Range("E38").Value = '  +1.37%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.995'
$ws.Range("E39").Value = '  +0.01%  '

$ws.Range("E40").Value = '  -2.93%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '44.81'
$ws.Range("E43").Value = '  -2.76%  '

$ws.Range("E44").Value = '  -1.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.300'
$ws.Range("E45").Value = '  +0.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.22'
$ws.Range("E46").Value = '  +6.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '151.78'
$ws.Range("E47").Value = '  +1.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.40'
$ws.Range("E48").Value = '  +12.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.36'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.85'
$ws.Range("E50").Value = '  +1.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '390.91'
$ws.Range("E51").Value = '  -0.37%  '
